$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '25.847.76'
$ws.Range("E2").Value = '  -0.12%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.641.19'
$ws.Range("E3").Value = '  +0.16%  '
$ws.Range("E4").Value = '  -0.29%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '215.83'
$ws.Range("E5").Value = '  +0.05%  '
$ws.Range("E6").Value = '  -0.60%  '
$ws.Range("E7").Value = '  -0.24%  '
$ws.Range("E8").Value = '  -0.73%  '
$ws.Range("E9").Value = '  -0.84%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.86'
$ws.Range("E10").Value = '  -1.85%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0794'
$ws.Range("E11").Value = '  +1.46%  '
$ws.Range("E12").Value = '  +0.42%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.868.02'
$ws.Range("E13").Value = '  +0.19%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.639.08'
$ws.Range("E14").Value = '  -0.39%  '
$ws.Range("E15").Value = '  +0.00%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0₃0770'
$ws.Range("E16").Value = '  +0.25%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '63.28'
$ws.Range("E17").Value = '  -0.04%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '25.867.12'
$ws.Range("E18").Value = '  -0.03%  '
$ws.Range("E19").Value = '  -0.28%  '
$ws.Range("E20").Value = '  +2.37%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '193.39'
$ws.Range("E21").Value = '  -0.50%  '
$ws.Range("E22").Value = '  +0.69%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.39'
$ws.Range("E23").Value = '  +2.85%  '
$ws.Range("E24").Value = '  +4.04%  '
$ws.Range("E25").Value = '  -0.24%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '143.02'
$ws.Range("E26").Value = '  +3.38%  '
$ws.Range("E27").Value = '  +0.33%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.00'
$ws.Range("E28").Value = '  +2.16%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.60'
$ws.Range("E29").Value = '  +0.25%  '
$ws.Range("E30").Value = '  -0.21%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0498'
$ws.Range("E31").Value = '  +0.66%  '
$ws.Range("E32").Value = '  +1.54%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.26'
$ws.Range("E33").Value = '  +0.32%  '
$ws.Range("E34").Value = '  +0.67%  '
$ws.Range("E35").Value = '  -0.41%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.913'
$ws.Range("E36").Value = '  +0.50%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.133.61'
$ws.Range("E37").Value = '  +0.62%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.53'
$ws.Range("E38").Value = '  -2.24%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.548'
$ws.Range("E39").Value = '  -0.62%  '
$ws.Range("E40").Value = '  -0.07%  '
$ws.Range("E41").Value = '  -0.06%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.57'
$ws.Range("E42").Value = '  +1.31%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '100.60'
$ws.Range("E43").Value = '  +1.15%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.808'
$ws.Range("E44").Value = '  +1.09%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.777.28'
$ws.Range("E45").Value = '  +0.11%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0₆0113'
$ws.Range("E46").Value = '  -0.14%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '55.59'
$ws.Range("E47").Value = '  +0.10%  '
$ws.Range("E48").Value = '  -1.40%  '

# Rows 49-51: reorder (RenderToken <-> Cronos swap) and replace Algorand with SynthetixNetwork
$ws.Range("B49").Value = 'RenderToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.44'
$ws.Range("E49").Value = '  +5.57%  '
$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0504'
$ws.Range("E50").Value = '  -0.17%  '
$ws.Range("B51").Value = 'SynthetixNetwork'
$ws.Range("C51").Value = 'https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.33'
$ws.Range("E51").Value = '  +3.63%  '
